# Update the three API endpoint URLs shown in the ANSWER column (column C)
# for the GA2.9 (FastAPI), GA3.7 (similarity), GA3.8 (execute) and GA4.3
# (outline) rows: the ports used in the example URLs were changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = "http://127.0.0.1:9000/api"
$ws.Range("C36").Value = "http://127.0.0.1:10000/similarity"
$ws.Range("C37").Value = "http://127.0.0.1:11000/execute"
$ws.Range("C41").Value = "http://127.0.0.1:12000/api/outline"

# Restore the view/selection state that was saved with the workbook
# (scrolled up a bit, with B26 as the active cell instead of B53).
$ws.Range("A18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B26").Select() | Out-Null
